# #327 Ajout des profils d'acces
# 1) Bump the "Date" metadata value.
# 2) Swap the two mapping columns on the Elements sheet so that the new
#    "Mapping: Specification metier vers l'extension ROR LocationSupportedCapacity"
#    column comes before the existing "Mapping: RIM Mapping" column (was AL, AK -> now AK, AL).

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 (Date row) ------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2) Elements sheet: swap columns AK (37) and AL (38) --------------------
$els = $wb.Worksheets.Item("Elements")

$lastRow = 61
$colAK = 37
$colAL = 38

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $els.Cells.Item($r, $colAK)
    $alCell = $els.Cells.Item($r, $colAL)

    $akText = $akCell.Text
    $alText = $alCell.Text

    $akCell.Value = $alText
    $alCell.Value = $akText
}

# Column widths follow the content: the wide "Mapping: Specification metier..."
# column is now AK (37) and the narrower "Mapping: RIM Mapping" column is AL (38).
$els.Columns.Item($colAK).ColumnWidth = 83.625
$els.Columns.Item($colAL).ColumnWidth = 24.98046875
